$d = $word.ActiveDocument

# Locate the "C1" text in the "Sección: C1" line (avoid hard-coded offsets).
$findRange = $d.Content
$found = $findRange.Find.Execute("C1", $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'C1' in the document"
}
$c1Start = $findRange.Start
$c1End = $findRange.End

# Step 1: change the text "C1" to "C2" (keeps it as a single run for now --
# this also coalesces it with the surrounding space runs into one run).
$c1Range = $d.Range($c1Start, $c1End)
$c1Range.Text = "C2"

$cEnd = $c1Start + 1      # position right after "C"
$twoEnd = $c1Start + 2    # position right after "2"

# Step 2: re-split the merged run back into separate "C" and "2" runs (and
# leave the surrounding space runs on their own) by toggling a character
# property on/off on each sub-range -- this forces the engine to materialize
# a distinct run at that boundary without re-merging everything the way a
# text insert/delete does.
$cRange = $d.Range($c1Start, $cEnd)
$cRange.Bold = 1
$cRange.Bold = 0

$twoRange = $d.Range($cEnd, $twoEnd)
$twoRange.Bold = 1
$twoRange.Bold = 0

# Step 3: move the "_GoBack" bookmark (Word's "last edit location" marker) to
# sit right after the new "2" run. Bookmarks.Add relocates a bookmark that
# already exists under that name instead of duplicating it, so this also
# removes the stale "_GoBack" bookmark that used to sit near the page break
# further down in the document.
$goBackRange = $d.Range($twoEnd, $twoEnd)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
